# Updates the "cryptos" price/volume table with freshly scraped values.
# Note: some Price (column D) values look like plain decimal numbers
# (e.g. "574.30"). Assigning such a string via .Value would make Excel
# auto-convert it to a numeric value (losing the trailing zero / exact
# text). To keep these as literal text - matching how the sheet already
# stores every other price as text - a leading apostrophe (text-prefix)
# is included for those particular values; Excel strips it from the
# stored value and only remembers that the cell was entered as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.737.84'
$ws.Range('E2').Value = '  +0.63%  '

$ws.Range('D3').Value = '3.384.63'
$ws.Range('E3').Value = '  +0.31%  '

$ws.Range('E4').Value = '  +0.37%  '

$ws.Range('D5').Value = '''574.30'
$ws.Range('E5').Value = '  -0.07%  '

$ws.Range('D6').Value = '''138.46'
$ws.Range('E6').Value = '  +0.52%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').Value = '3.383.33'
$ws.Range('E8').Value = '  +0.26%  '

$ws.Range('D9').Value = '''0.474'
$ws.Range('E9').Value = '  -0.64%  '

$ws.Range('D10').Value = '''7.47'
$ws.Range('E10').Value = '  -1.93%  '

$ws.Range('D11').Value = '''0.124'
$ws.Range('E11').Value = '  -0.31%  '

$ws.Range('D12').Value = '''0.392'
$ws.Range('E12').Value = '  +0.23%  '

$ws.Range('D13').Value = '3.960.63'
$ws.Range('E13').Value = '  +0.52%  '

$ws.Range('D14').Value = '''0.124'
$ws.Range('E14').Value = '  +1.79%  '

$ws.Range('D15').Value = '''0.0000175'
$ws.Range('E15').Value = '  -0.93%  '

$ws.Range('D16').Value = '''26.00'
$ws.Range('E16').Value = '  +2.41%  '

$ws.Range('D17').Value = '3.385.65'
$ws.Range('E17').Value = '  +0.38%  '

$ws.Range('D18').Value = '61.874.44'
$ws.Range('E18').Value = '  +0.87%  '

$ws.Range('E19').Value = '  +1.02%  '

$ws.Range('D20').Value = '''13.97'
$ws.Range('E20').Value = '  +0.11%  '

$ws.Range('D21').Value = '''9.41'
$ws.Range('E21').Value = '  +0.48%  '

$ws.Range('D22').Value = '''377.93'
$ws.Range('E22').Value = '  -1.09%  '

$ws.Range('D23').Value = '''0.557'
$ws.Range('E23').Value = '  -2.15%  '

$ws.Range('D24').Value = '3.521.95'
$ws.Range('E24').Value = '  +0.57%  '

$ws.Range('E25').Value = '  -0.07%  '

$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '''71.36'
$ws.Range('E26').Value = '  +1.22%  '

$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '''0.0000124'
$ws.Range('E27').Value = '  +0.88%  '

$ws.Range('D28').Value = '''1.80'
$ws.Range('E28').Value = '  +10.93%  '

$ws.Range('D29').Value = '''7.64'
$ws.Range('E29').Value = '  -2.54%  '

$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.57%  '

$ws.Range('D31').Value = '''0.164'
$ws.Range('E31').Value = '  +4.04%  '

$ws.Range('D32').Value = '''8.23'
$ws.Range('E32').Value = '  -0.46%  '

$ws.Range('D33').Value = '''2.16'
$ws.Range('E33').Value = '  +1.37%  '

$ws.Range('E34').Value = '  +0.09%  '

$ws.Range('D35').Value = '''23.64'
$ws.Range('E35').Value = '  +0.68%  '

$ws.Range('D36').Value = '''5.21'
$ws.Range('E36').Value = '  -5.49%  '

$ws.Range('D37').Value = '''6.84'
$ws.Range('E37').Value = '  -3.05%  '

$ws.Range('D38').Value = '''1.54'
$ws.Range('E38').Value = '  +0.02%  '

$ws.Range('D39').Value = '''165.04'
$ws.Range('E39').Value = '  +2.49%  '

$ws.Range('D40').Value = '''0.0771'
$ws.Range('E40').Value = '  -2.46%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.21%  '

$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''1.73'
$ws.Range('E42').Value = '  -0.14%  '

$ws.Range('D43').Value = '''0.775'
$ws.Range('E43').Value = '  +1.38%  '

$ws.Range('D44').Value = '''1.21'
$ws.Range('E44').Value = '  -0.52%  '

$ws.Range('D45').Value = '''41.64'
$ws.Range('E45').Value = '  +0.42%  '

$ws.Range('D46').Value = '''4.39'
$ws.Range('E46').Value = '  -1.04%  '

$ws.Range('D47').Value = '''24.15'
$ws.Range('E47').Value = '  +3.75%  '

$ws.Range('D48').Value = '''6.87'
$ws.Range('E48').Value = '  -1.44%  '

$ws.Range('E49').Value = '  +1.23%  '

$ws.Range('D50').Value = '2.387.23'
$ws.Range('E50').Value = '  +2.09%  '

$ws.Range('E51').Value = '  +0.49%  '

